# Insert a new data row at row 163 (pushing the existing rows 163:250 down
# to 164:251) and populate the new row with the latest price-report entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(163).Insert()

$ws.Cells.Item(163, 1).Value = 7
$ws.Cells.Item(163, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(163, 3).Value = "Ñuble"
$ws.Cells.Item(163, 4).Value = 45001
$ws.Cells.Item(163, 5).Value = 16
$ws.Cells.Item(163, 6).Value = "Fruta"
$ws.Cells.Item(163, 7).Value = 100104
$ws.Cells.Item(163, 8).Value = "Frutos de pepita"
$ws.Cells.Item(163, 9).Value = 100104005
$ws.Cells.Item(163, 10).Value = "Pera"
$ws.Cells.Item(163, 11).Value = "Packham's Triumph"
$ws.Cells.Item(163, 12).Value = "Primera"
$ws.Cells.Item(163, 13).Value = 50
$ws.Cells.Item(163, 14).Value = 10000
$ws.Cells.Item(163, 15).Value = 10000
$ws.Cells.Item(163, 16).Value = 10000
$ws.Cells.Item(163, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(163, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(163, 19).Value = 556
$ws.Cells.Item(163, 20).Value = 18
